$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B14 value (semana 13) from 478 to 490
$ws.Range("B14").Value = 490

# Add new row 15 for semana 14, casos 734
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 734
